$d = $word.ActiveDocument

# --- Programa (Portuguese) paragraph ---
# Insert a line break between "...hiperbólicas. " and "Matrizes..."
$d.Content.Find.Execute(
    "hiperbólicas. Matrizes, determinantes e sistemas lineares",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "hiperbólicas. ^lMatrizes, determinantes e sistemas lineares",
    2)

# Insert a line break between "...sistemas lineares" and "Modelagem:..."
$d.Content.Find.Execute(
    "sistemas linearesModelagem: Áreas",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "sistemas lineares^lModelagem: Áreas",
    2)

# --- Programa (English, italic) paragraph ---
# Insert a line break between "...inequalities." and "Real Functions:..."
$d.Content.Find.Execute(
    "inequalities.Real Functions: Definition of functions",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "inequalities.^lReal Functions: Definition of functions",
    2)

# Insert a line break between "...hyperbolic functions." and "Matrix,..."
$d.Content.Find.Execute(
    "hyperbolic functions.Matrix, determinant and linear systems.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "hyperbolic functions.^lMatrix, determinant and linear systems.",
    2)

# Insert a line break between "...linear systems." and "Modeling:..."
$d.Content.Find.Execute(
    "linear systems.Modeling: Areas, volume, cost, population models.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "linear systems.^lModeling: Areas, volume, cost, population models.",
    2)

# --- Bibliografia paragraph ---
# Insert a blank line (double break) between each reference entry.
$d.Content.Find.Execute(
    "v.1.ANTON, Howard.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "v.1.^l^lANTON, Howard.",
    2)

$d.Content.Find.Execute(
    "Bookman, 2007.THOMAS, George B.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Bookman, 2007.^l^lTHOMAS, George B.",
    2)

$d.Content.Find.Execute(
    "v.1,FLEMMING, Diva M.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "v.1,^l^lFLEMMING, Diva M.",
    2)
